# Update the VIN value used in column A (rows 2-5 all share the same
# string "ZZZZN3DD&E" -> "FFFKN3DD&E").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A5").Value = "FFFKN3DD&E"

# Move/save the active selection as it was when the workbook was last saved.
$ws.Range("A8").Select() | Out-Null
